$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the formatting (border/wrap/vertical-top style already used by the
# header row) down through the two new data rows before filling in values,
# so every new cell picks up the same style index as the rest of the table.
$ws.Range("A1:K1").Copy()
$ws.Range("A2:K3").PasteSpecial(-4122)

# Row 2: Botswana AMR Gram Stain Observation
$ws.Range("A2").Value = "botswana-amr-gram-stain-observation"
$ws.Range("B2").Value = "Botswana AMR Gram Stain Observation"
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "LOINC#664-3"
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = "dateTime"
$ws.Range("H2").Value = "CodeableConcept"
$ws.Range("I2").Value = "optional"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""

# Row 3: Botswana AMR Organism Identification Observation
$ws.Range("A3").Value = "botswana-amr-organism-observation"
$ws.Range("B3").Value = "Botswana AMR Organism Identification Observation"
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "LOINC#634-6"
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = "dateTime"
$ws.Range("H3").Value = "CodeableConcept"
$ws.Range("I3").Value = "optional"
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = ""
